$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 132, shifting existing rows 132..304 down to 133..305
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with a new weekly price observation
$ws.Range("A132").Value = 3
$ws.Range("B132").Value = "Femacal de La Calera"
$ws.Range("C132").Value = "Coquimbo"
$ws.Range("D132").Value = 44671
$ws.Range("E132").Value = 5
$ws.Range("F132").Value = 100112039
$ws.Range("G132").Value = "Ciboulette"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 150
$ws.Range("K132").Value = 1500
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = 1500
$ws.Range("N132").Value = "$/docena de atados"
$ws.Range("O132").Value = "Provincia de Quillota"
$ws.Range("P132").Value = 500
$ws.Range("Q132").Value = 3
$ws.Range("R132").Value = "Hortaliza"

# Match the date-number style (s="2") used by the rest of column D
$ws.Range("D132").NumberFormat = $ws.Range("D133").NumberFormat
